$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity for the Kingbright SC10-21SRWA (7-Segment display) row (row 6)
$ws.Range("D6").Value = 12

# Update the 5X Board Quantity formula for row 6 to reflect sourcing changes
$ws.Range("E6").Formula = "=8*_xlfn.CEILING.MATH(20/3)+4*5"

# Slightly widen column G (Total Cost 5X) to fit the new values
# (Excel COM stores column width on a whole-pixel grid, so 12.5 "chars" is
# the closest settable value that rounds to the target 13.35 stored width)
$ws.Columns.Item(7).ColumnWidth = 12.5

# Move the active selection from B8 to F8
$ws.Range("F8").Select()
